# Word COM-interop script: swap the "Virtual Reality" essay for the
# "Ecosystem equilibrium" essay (title, byline, e-mail, body, summary),
# and append a trailing empty paragraph, as described by the source diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title (paragraph 1)
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute(
    "Virtual Reality: An Immersive Paradigm Shift", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Understanding the Dynamic Equilibrium of Ecosystems: A Balanced Orchestra of Life",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Byline / author name (paragraph 2): "Varun Wadhwani" -> "Dr. Alexander Westwood"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute(
    "Varun Wadhwani", $true, $false, $false, $false, $false, $true, 1,
    $false, "Dr. Alexander Westwood", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) E-mail address (paragraph 3): varun.wadhwani@emailworld.com
#                                -> westwood.a@eduinstitute.org
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("varun", $true, $false, $false, $false, $false, $true, 1,
    $false, "westwood", 2) | Out-Null

$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("wadhwani@emailworld", $true, $false, $false, $false,
    $false, $true, 1, $false, "a@eduinstitute", 2) | Out-Null

$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("com", $true, $false, $false, $false, $false, $true, 1,
    $false, "org", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Body paragraph (paragraph 5) - three sentence-groups separated by
#    manual line breaks. Replace each group's full text in one shot.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(5).Range

$old5a = "Virtual Reality(VR), a captivating technology, transports users into a simulated environment, empowering them to interact with digital creations. Its transformative properties have sparked a paradigm shift across various disciplines, from gaming and entertainment to education, healthcare, and workplace simulations. As VR's applications continue to expand, let's delve into its profound impact and explore the boundless possibilities it holds for shaping the future."
$new5a = "The intricate harmony of ecosystems lies in the delicate balance maintained between organisms and their environment. Life thrives within this dynamic equilibrium, where innumerable interactions weave a complex tapestry of interdependence. Every organism plays a unique role in maintaining this delicate balance, contributing to the overall stability and resilience of the ecosystem. To comprehend the intricate dance of life, we must delve into the fundamental principles governing these interactions."
$p5.Find.Execute($old5a, $true, $false, $false, $false, $false, $true, 1,
    $false, $new5a, 2) | Out-Null

$p5 = $d.Paragraphs(5).Range
$old5b = "VR's immersive nature has revolutionized gaming and entertainment, providing an unparalleled level of engagement and escapism. Players can now step into virtual worlds, embodying characters and experiencing adventures like never before. This immersive experience has also found its way into other fields, such as education and training, where simulations can provide realistic and engaging scenarios for students and professionals alike."
$new5b = "Understanding the roles of individual organisms within an ecosystem is crucial. Each species occupies a specific ecological niche, playing a distinct role in energy flow and nutrient cycling. This interconnectedness forms a web of relationships that shape the dynamics of the ecosystem. Changes in one species can ripple through the entire system, triggering a cascade of ecological responses. By studying these interactions, scientists can unravel the intricate mechanisms that maintain equilibrium."
$p5.Find.Execute($old5b, $true, $false, $false, $false, $false, $true, 1,
    $false, $new5b, 2) | Out-Null

$p5 = $d.Paragraphs(5).Range
$old5c = "Beyond entertainment and education, VR is making significant strides in healthcare. It offers immersive therapies for conditions like PTSD and phobias, allowing patients to confront their fears in a controlled virtual environment. Moreover, VR is instrumental in surgical training, enabling surgeons to practice complex procedures in a risk-free environment, leading to improved surgical outcomes."
$new5c = "Biodiversity, the vast array of species within an ecosystem, is crucial for maintaining ecological balance. A rich diversity of species enhances the resilience of ecosystems, allowing them to adapt to environmental changes. By providing a variety of habitats and resources, biodiversity ensures the survival of a wide range of organisms. Preserving biodiversity is therefore essential for the long-term stability and productivity of ecosystems."
$p5.Find.Execute($old5c, $true, $false, $false, $false, $false, $true, 1,
    $false, $new5c, 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Summary paragraph (paragraph 7)
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs(7).Range
$old7 = "Virtual Reality has emerged as a pivotal technology, ushering in a new era of immersive experiences across various domains. Its applications have soared, ranging from gaming and entertainment to education, healthcare, and corporate training. VR's ability to transport users into digital environments has enabled novel and engaging ways of learning, healing, and simulating complex scenarios. As technology continues to advance, VR's impact is poised to grow exponentially, reshaping industries and transforming the way we learn, heal, and experience the world around us."
$new7 = "Ecosystems are complex systems in which organisms and their environment interact, creating a dynamic balance. The roles of individual species, their interconnectedness, and biodiversity are key factors in maintaining this balance. Understanding these interactions allows us to appreciate the delicate harmony of life and the importance of preserving biodiversity. By studying ecosystems, we gain insights into the intricate dance of life, unraveling the secrets of maintaining a healthy and balanced environment."
$p7.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1,
    $false, $new7, 2) | Out-Null

# ---------------------------------------------------------------------
# 6) Append a new, trailing empty paragraph at the end of the document.
# ---------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null
